# Update countries & provincias Spain
# Applies the 16-Oct-2020 data refresh (13:15 -> 14:32) to the "Pais" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update the "last updated" timestamp banner (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 16 de Octubre de 2020 a las 14:32"

# --- Estados Unidos (row 4) ---
$ws.Range("B4").Value = 8219831
$ws.Range("C4").Value = 3516
$ws.Range("D4").Value = 5320747
$ws.Range("E4").Value = 2676303
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 64
$ws.Range("H4").Value = 222781

# --- Paises Bajos (row 29) ---
$ws.Range("B29").Value = 211938
$ws.Range("C29").Value = 7984
$ws.Range("G29").Value = 16
$ws.Range("H29").Value = 6708

# --- Suecia (row 47) ---
$ws.Range("B47").Value = 103200
$ws.Range("G47").Value = 6
$ws.Range("H47").Value = 5918

# --- Azerbaiyan (row 73) ---
$ws.Range("B73").Value = 43789
$ws.Range("C73").Value = 509
$ws.Range("D73").Value = 39800
$ws.Range("E73").Value = 3368
$ws.Range("G73").Value = 2
$ws.Range("H73").Value = 621

# --- Dinamarca (row 79) ---
$ws.Range("B79").Value = 34441
$ws.Range("C79").Value = 418
$ws.Range("D79").Value = 28551
$ws.Range("E79").Value = 5213

# --- Finlandia (row 102) ---
$ws.Range("E102").Value = 3682
$ws.Range("G102").Value = 1
$ws.Range("H102").Value = 351

# --- Rows 110-112: Uganda overtakes Luxemburgo & Guayana Francesa ---
# Row 110 becomes Uganda (updated figures)
$ws.Range("A110").Value = "Uganda"
$ws.Range("B110").Value = 10334
$ws.Range("C110").Value = 217
$ws.Range("D110").Value = 6901
$ws.Range("E110").Value = 3337
$ws.Range("F110").Value = 0
$ws.Range("G110").Value = 0
$ws.Range("H110").Value = 96

# Row 111 becomes Luxemburgo (old row-110 figures)
$ws.Range("A111").Value = "Luxemburgo"
$ws.Range("B111").Value = 10244
$ws.Range("C111").Value = 0
$ws.Range("D111").Value = 8384
$ws.Range("E111").Value = 1727
$ws.Range("F111").Value = 0
$ws.Range("G111").Value = 0
$ws.Range("H111").Value = 133

# Row 112 becomes Guayana Francesa (old row-111 figures)
$ws.Range("A112").Value = "Guayana Francesa"
$ws.Range("B112").Value = 10233
$ws.Range("C112").Value = 0
$ws.Range("D112").Value = 9894
$ws.Range("E112").Value = 270
$ws.Range("F112").Value = 0
$ws.Range("G112").Value = 0
$ws.Range("H112").Value = 69

# --- Sri Lanka (row 128) ---
$ws.Range("B128").Value = 5305
$ws.Range("C128").Value = 61
$ws.Range("E128").Value = 1907

# --- Islandia (row 142) ---
$ws.Range("E142").Value = 1205
$ws.Range("G142").Value = 1
$ws.Range("H142").Value = 11

Write-Host "Applied paises.xlsx data refresh (16 Oct 2020, 14:32)"
